$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2486.96
$ws.Range("I40").Value = 2489.0908
$ws.Range("J40").Value = 2485.2856
$ws.Range("K40").Value = 2489.0908
$ws.Range("L40").Value = 2485.2856
$ws.Range("M40").Value = -2314.0908
$ws.Range("N40").Value = -2835.2856

$ws.Range("H64").Value = 3633.2354
$ws.Range("I64").Value = 3340.3572
$ws.Range("K64").Value = 3340.3572
$ws.Range("M64").Value = -3092.3572

$ws.Range("H67").Value = 3633.2354
$ws.Range("I67").Value = 3340.3572
$ws.Range("K67").Value = 3340.3572
$ws.Range("M67").Value = -2482.3572

$ws.Range("H76").Value = 3381.3416
$ws.Range("I76").Value = 2711.0688
$ws.Range("K76").Value = 2711.0688
$ws.Range("M76").Value = -2396.0688

$ws.Range("H79").Value = 3381.3416
$ws.Range("I79").Value = 2711.0688
$ws.Range("K79").Value = 2711.0688
$ws.Range("M79").Value = -1619.0688

$ws.Range("H87").Value = 13053.356
$ws.Range("J87").Value = 13053.356
$ws.Range("L87").Value = 13053.356
$ws.Range("N87").Value = -15549.356

$ws.Range("H90").Value = 13053.356
$ws.Range("J90").Value = 13053.356
$ws.Range("L90").Value = 39160.068
$ws.Range("N90").Value = -51640.068

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2885.54
$ws.Range("I32").Value = 2885.54
$ws.Range("K32").Value = 2885.54
$ws.Range("M32").Value = -2598.54

$ws.Range("H122").Value = 2214.2
$ws.Range("I122").Value = 1955.5454
$ws.Range("K122").Value = 5866.6362
$ws.Range("M122").Value = -3416.6362

$ws.Range("H132").Value = 2290.8647
$ws.Range("I132").Value = 2112.8948
$ws.Range("J132").Value = 2478.7222
$ws.Range("K132").Value = 6338.6844
$ws.Range("L132").Value = 7436.1666
$ws.Range("M132").Value = -3808.6844
$ws.Range("N132").Value = -12496.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2818.3845
$ws.Range("I105").Value = 2846.9
$ws.Range("J105").Value = 2723.3333
$ws.Range("K105").Value = 2846.9
$ws.Range("L105").Value = 2723.3333
$ws.Range("M105").Value = -1099.9
$ws.Range("N105").Value = -6217.3333

$ws.Range("H140").Value = 49870
$ws.Range("J140").Value = 49870
$ws.Range("L140").Value = 49870
$ws.Range("N140").Value = -60230

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6019.636
$ws.Range("I62").Value = 4415.7144
$ws.Range("J62").Value = 8826.5
$ws.Range("K62").Value = 4415.7144
$ws.Range("L62").Value = 8826.5
$ws.Range("M62").Value = -3791.7144
$ws.Range("N62").Value = -10074.5

$ws.Range("H65").Value = 6019.636
$ws.Range("I65").Value = 4415.7144
$ws.Range("J65").Value = 8826.5
$ws.Range("K65").Value = 22078.572
$ws.Range("L65").Value = 44132.5
$ws.Range("M65").Value = -18958.572
$ws.Range("N65").Value = -50372.5

$ws.Range("H140").Value = 80335.60000000001
$ws.Range("J140").Value = 80335.60000000001
$ws.Range("L140").Value = 80335.60000000001
$ws.Range("N140").Value = -90695.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 61.72222
$ws.Range("I12").Value = 100.666664
$ws.Range("J12").Value = 53.933334
$ws.Range("K12").Value = 301.999992
$ws.Range("L12").Value = 161.800002
$ws.Range("M12").Value = -128.999992
$ws.Range("N12").Value = -507.800002

$ws.Range("H19").Value = 2208.182
$ws.Range("J19").Value = 2400
$ws.Range("L19").Value = 7200
$ws.Range("N19").Value = -7548

$ws.Range("H37").Value = 97550
$ws.Range("J37").Value = 97550
$ws.Range("L37").Value = 292650
$ws.Range("N37").Value = -292874

$ws.Range("H80").Value = 3099.6667
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3178.2144
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 9534.643199999999
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -11406.6432

$ws.Range("H83").Value = 3099.6667
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3178.2144
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 28603.9296
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -37963.9296

$ws.Range("H131").Value = 834.3
$ws.Range("I131").Value = 384
$ws.Range("J131").Value = 858
$ws.Range("K131").Value = 1152
$ws.Range("L131").Value = 2574
$ws.Range("M131").Value = 3888
$ws.Range("N131").Value = -12654

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5706.2354
$ws.Range("I70").Value = 4779.1875
$ws.Range("K70").Value = 4779.1875
$ws.Range("M70").Value = -4509.1875

$ws.Range("H73").Value = 5706.2354
$ws.Range("I73").Value = 4779.1875
$ws.Range("K73").Value = 4779.1875
$ws.Range("M73").Value = -3843.1875

$ws.Range("H122").Value = 2547.9736
$ws.Range("I122").Value = 1991.4642
$ws.Range("J122").Value = 4106.2
$ws.Range("K122").Value = 5974.392599999999
$ws.Range("L122").Value = 12318.6
$ws.Range("M122").Value = -3524.392599999999
$ws.Range("N122").Value = -17218.6

$ws.Range("H132").Value = 2396.2424
$ws.Range("I132").Value = 2024.5
$ws.Range("J132").Value = 2968.1538
$ws.Range("K132").Value = 6073.5
$ws.Range("L132").Value = 8904.4614
$ws.Range("M132").Value = -3543.5
$ws.Range("N132").Value = -13964.4614

$ws.Range("H135").Value = 39037.777
$ws.Range("J135").Value = 39037.777
$ws.Range("L135").Value = 39037.777
$ws.Range("N135").Value = -49177.777

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41288.92
$ws.Range("I7").Value = 52348.15
$ws.Range("K7").Value = 52348.15
$ws.Range("M7").Value = -52236.15

$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 700
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1290

$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 700
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 700
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -914

$ws.Range("H126").Value = 41288.92
$ws.Range("I126").Value = 52348.15
$ws.Range("K126").Value = 157044.45
$ws.Range("M126").Value = -154574.45

$ws.Range("H132").Value = 9666.471
$ws.Range("I132").Value = 8319.941999999999
$ws.Range("J132").Value = 12612
$ws.Range("K132").Value = 24959.826
$ws.Range("L132").Value = 37836
$ws.Range("M132").Value = -22429.826
$ws.Range("N132").Value = -42896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1244.6666
$ws.Range("I132").Value = 902.4474
$ws.Range("J132").Value = 2057.4375
$ws.Range("K132").Value = 2707.3422
$ws.Range("L132").Value = 6172.3125
$ws.Range("M132").Value = -177.3422
$ws.Range("N132").Value = -11232.3125
